# Pavani Gandepalli Status on 21/Jan/2021
# Fill in the new status row (row 17) that documents the 21/1/2021 entry:
#   A17 -> date, B17 -> the day's notes, D17 -> attached file name.
# The shared-string table must grow in the same order as the source edit:
# date (28), filename (29), notes (30) - so we write A17, then D17, then B17.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A17").Value = "21/1/2021"
$ws.Range("D17").Value = "BirthdayCake.txt"

$notes = "1. Learnt what is pointer, how to declare, initialize and access pointer variables" + "`n" + `
         "2. What is NULL pointer" + "`n" + `
         "3. Pointer Arthematic " + "`n" + `
         "4. Strings, How to declare, initialize and access strings" + "`n" + `
         "5. Using %s, scanf, gets, puts, fgets for string input and optput" + "`n" + `
         "6. Uderstood predefined functions in strings" + "`n" + `
         "7. 2-D strings , memory allocation of strings " + "`n" + `
         "8. Completed one hacker rank program" + "`n" + `
         "9. Attended session by Srinivasa"

$ws.Range("B17").Value = $notes

# Match the formatting used by the row above (row 16), which has the same
# A/B/C/D layout: date style, wrapped notes style, blank, filename style.
$ws.Range("A17").Style = $ws.Range("A16").Style
$ws.Range("B17").Style = $ws.Range("B16").Style
$ws.Range("D17").Style = $ws.Range("D16").Style

$ws.Rows.Item(17).RowHeight = 180

# Update the view: scrolled down one row further, with C17 now selected.
$ws.Range("C17").Select()
